$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no numeric auto-coercion) for Price column cells
# that would otherwise look like plain numbers to Excel.

$ws.Range("D2").Value = '29.284.62'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.843.80'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.04'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07439'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2940'
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.91'
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = '1.813.95'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.008'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6715'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.84'
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.151'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '29.279.84'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008319'
$ws.Range("E18").Value = '  +2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.65'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.52'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.162'
$ws.Range("E22").Value = '  -2.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '161.01'
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.700'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1405'
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.02'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.515'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.157'
$ws.Range("E29").Value = '  -2.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.066'
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05296'
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7521'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.137'
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.680'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").Value = '1.319.10'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01806'
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9213'
$ws.Range("E40").Value = '  -1.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.977'
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.08334'
$ws.Range("E42").Value = '  +13.36%  '
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.16'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.989.16'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5167'
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.775'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000121'
$ws.Range("E48").Value = '  -6.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.82'
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.145'
$ws.Range("E50").Value = '  -3.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05949'
$ws.Range("E51").Value = '  +0.08%  '
